$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# ---------------------------------------------------------------------------
# 1) Rows 12 & 18 (raw "input" rows) get four new sprint/semana columns
#    (M:P). Copy the formatting from the last populated cell in each row
#    (column L, style s="5" - red Arial) and then write the new numbers.
#    Everything downstream (rows 10,11,13,14,15,16,17) is driven by shared
#    formulas already present in the template, so they recalc automatically.
# ---------------------------------------------------------------------------
$ws.Range("L12").Copy()
$ws.Range("M12:P12").PasteSpecial(-4122)
$ws.Range("M12").Value = 87.7
$ws.Range("N12").Value = 88.7
$ws.Range("O12").Value = 88.7
$ws.Range("P12").Value = 104.9

$ws.Range("L18").Copy()
$ws.Range("M18:P18").PasteSpecial(-4122)
$ws.Range("M18").Value = 85.7
$ws.Range("N18").Value = 87.2
$ws.Range("O18").Value = 87.2
$ws.Range("P18").Value = 108.2

# ---------------------------------------------------------------------------
# 2) Small side table in columns R:T (rows 14-20) - labels + minutes + a
#    total-hours formula in T20.
# ---------------------------------------------------------------------------
$ws.Range("G5").Copy()
$ws.Range("R14").PasteSpecial(-4122)
$ws.Range("R15").PasteSpecial(-4122)
$ws.Range("R16").PasteSpecial(-4122)
$ws.Range("R17").PasteSpecial(-4122)
$ws.Range("R18").PasteSpecial(-4122)
$ws.Range("R19").PasteSpecial(-4122)
$ws.Range("R20").PasteSpecial(-4122)
$ws.Range("S14").PasteSpecial(-4122)
$ws.Range("S15").PasteSpecial(-4122)
$ws.Range("S16").PasteSpecial(-4122)
$ws.Range("S17").PasteSpecial(-4122)
$ws.Range("S18").PasteSpecial(-4122)
$ws.Range("S19").PasteSpecial(-4122)
$ws.Range("S20").PasteSpecial(-4122)
$ws.Range("M19").PasteSpecial(-4122)
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("W15").PasteSpecial(-4122)

$ws.Range("R14").Value = "P IA"
$ws.Range("S14").Value = 360.0
$ws.Range("R15").Value = "BD"
$ws.Range("S15").Value = 345.0
$ws.Range("R16").Value = "LFG"
$ws.Range("S16").Value = 30.0
$ws.Range("R17").Value = "VE"
$ws.Range("S17").Value = 180.0
$ws.Range("R18").Value = "bf"
$ws.Range("S18").Value = 90.0
$ws.Range("R19").Value = "ct"
$ws.Range("S19").Value = 240.0
$ws.Range("R20").Value = "mv"
$ws.Range("S20").Value = 15.0

# Old standalone T15 / T20 placeholders are replaced by the table above.
$ws.Range("T15").Clear()
$ws.Range("T20").Clear()
$ws.Range("G5").Copy()
$ws.Range("T20").PasteSpecial(-4122)
$ws.Range("T20").Formula = "=SUM(S14:S20)/60"

# ---------------------------------------------------------------------------
# 3) Rows 23 & 24 (EV / AC helper rows) gain the same M:P columns as row 22,
#    mirroring rows 10 and 13 respectively. Style matches the existing I23:L23
#    / I24:L24 cells (s="6"), except P24 which keeps the plain s="1" style.
# ---------------------------------------------------------------------------
$ws.Range("I23").Copy()
$ws.Range("M23:P23").PasteSpecial(-4122)
$ws.Range("M23").Value = $ws.Range("M10").Value2
$ws.Range("N23").Value = $ws.Range("N10").Value2
$ws.Range("O23").Value = $ws.Range("O10").Value2
$ws.Range("P23").Value = $ws.Range("P10").Value2

$ws.Range("I24").Copy()
$ws.Range("M24:O24").PasteSpecial(-4122)
$ws.Range("M24").Value = $ws.Range("M13").Value2
$ws.Range("N24").Value = $ws.Range("N13").Value2
$ws.Range("O24").Value = $ws.Range("O13").Value2

$ws.Range("L22").Copy()
$ws.Range("P24").PasteSpecial(-4122)
$ws.Range("P24").Value = $ws.Range("P13").Value2

# ---------------------------------------------------------------------------
# 4) Rows 26 & 27 (SPI / CPI) - turn the existing E:L formulas into plain
#    cached values, then extend the row with M:P literal values (mirroring
#    rows 14 and 16 respectively).
# ---------------------------------------------------------------------------
foreach ($col in @("E","F","G","H","I","J","K","L")) {
    $cell26 = $ws.Range($col + "26")
    $v26 = $cell26.Value2
    $cell26.Value = $v26

    $cell27 = $ws.Range($col + "27")
    $v27 = $cell27.Value2
    $cell27.Value = $v27
}

$ws.Range("I23").Copy()
$ws.Range("M26:P26").PasteSpecial(-4122)
$ws.Range("M26").Value = $ws.Range("M14").Value2
$ws.Range("N26").Value = $ws.Range("N14").Value2
$ws.Range("O26").Value = $ws.Range("O14").Value2
$ws.Range("P26").Value = $ws.Range("P14").Value2

$ws.Range("I23").Copy()
$ws.Range("M27:P27").PasteSpecial(-4122)
$ws.Range("M27").Value = $ws.Range("M16").Value2
$ws.Range("N27").Value = $ws.Range("N16").Value2
$ws.Range("O27").Value = $ws.Range("O16").Value2
$ws.Range("P27").Value = $ws.Range("P16").Value2

# ---------------------------------------------------------------------------
# 5) Chart 2 (SPI/CPI line chart) now plots through column P instead of L.
# ---------------------------------------------------------------------------
$chart2 = $ws.ChartObjects(2).Chart
$chart2.SeriesCollection(1).Values = "=Planilha1!`$E`$26:`$P`$26"
$chart2.SeriesCollection(2).Values = "=Planilha1!`$E`$27:`$P`$27"
